$d = $word.ActiveDocument

# Merge the two runs "Graphical display for xrandr " + "and auto managment"
# into a single run with text "Graphical display for xrandr" (trailing space removed,
# "and auto managment" removed entirely).
$d.Content.Find.Execute("Graphical display for xrandr and auto managment", $true, $false, $false, $false, $false, $true, 1, $false, "Graphical display for xrandr", 2)

# Remove " autorandr" from the pacman install command.
$d.Content.Find.Execute("sudo pacman -S arandr autorandr", $true, $false, $false, $false, $false, $true, 1, $false, "sudo pacman -S arandr", 2)
